# Update the "想去人数" (want-to-go count) figures in column F across the
# relevant worksheets, matching the regenerated stats snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 35
$ws1.Range("F11").Value = 421
$ws1.Range("F14").Value = 897
$ws1.Range("F16").Value = 2022
$ws1.Range("F17").Value = 513
$ws1.Range("F18").Value = 8071
$ws1.Range("F19").Value = 663
$ws1.Range("F20").Value = 524
$ws1.Range("F21").Value = 68
$ws1.Range("F24").Value = 234

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5570
$ws3.Range("F3").Value = 409
$ws3.Range("F4").Value = 398

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5570
$ws4.Range("F4").Value = 409
$ws4.Range("F5").Value = 398
$ws4.Range("F11").Value = 35
$ws4.Range("F19").Value = 421
$ws4.Range("F24").Value = 897
$ws4.Range("F28").Value = 2022
$ws4.Range("F29").Value = 513
$ws4.Range("F30").Value = 8071
$ws4.Range("F33").Value = 663
$ws4.Range("F34").Value = 524
$ws4.Range("F35").Value = 68
$ws4.Range("F40").Value = 234
